# Auto-generated edit script applying scheduled-runner value refresh
# to Behemoth_Profits workbook (per-sheet leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17,J17,L17,N17
$ws.Cells.Item(17, 8).Value = 426.2
$ws.Cells.Item(17, 10).Value = 447.02127
$ws.Cells.Item(17, 12).Value = 1341.06381
$ws.Cells.Item(17, 14).Value = -1677.06381
# Row 33: H33,I33,K33,M33
$ws.Cells.Item(33, 8).Value = 1270.6786
$ws.Cells.Item(33, 9).Value = 1068.9166
$ws.Cells.Item(33, 11).Value = 1068.9166
$ws.Cells.Item(33, 13).Value = -839.9166
# Row 43: H43,I43,K43,M43
$ws.Cells.Item(43, 8).Value = 3265.7144
$ws.Cells.Item(43, 9).Value = 1899.0667
$ws.Cells.Item(43, 11).Value = 1899.0667
$ws.Cells.Item(43, 13).Value = -1830.0667
# Row 53: H53,I53,J53,K53,L53,M53,N53
$ws.Cells.Item(53, 8).Value = 321.54544
$ws.Cells.Item(53, 9).Value = 42.333332
$ws.Cells.Item(53, 10).Value = 426.25
$ws.Cells.Item(53, 11).Value = 42.333332
$ws.Cells.Item(53, 12).Value = 426.25
$ws.Cells.Item(53, 13).Value = 594.666668
$ws.Cells.Item(53, 14).Value = -1700.25
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Cells.Item(132, 8).Value = 6250.3335
$ws.Cells.Item(132, 9).Value = 6900.4
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 20701.2
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -18171.2
$ws.Cells.Item(132, 14).Value = -14060
# Row 137: H137,I137,K137,M137
$ws.Cells.Item(137, 8).Value = 4765.4287
$ws.Cells.Item(137, 9).Value = 1877.7142
$ws.Cells.Item(137, 11).Value = 5633.142599999999
$ws.Cells.Item(137, 13).Value = -3083.142599999999
# Row 138: H138,I138,K138,M138
$ws.Cells.Item(138, 8).Value = 2038.4681
$ws.Cells.Item(138, 9).Value = 1106.2222
$ws.Cells.Item(138, 11).Value = 3318.6666
$ws.Cells.Item(138, 13).Value = 1821.3334
# Row 141: H141,I141,K141,M141
$ws.Cells.Item(141, 8).Value = 3798.8
$ws.Cells.Item(141, 9).Value = 3798.8
$ws.Cells.Item(141, 11).Value = 11396.4
$ws.Cells.Item(141, 13).Value = -6216.400000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 45: H45,I45,K45,M45
$ws.Cells.Item(45, 8).Value = 3750.875
$ws.Cells.Item(45, 9).Value = 3715.2856
$ws.Cells.Item(45, 11).Value = 3715.2856
$ws.Cells.Item(45, 13).Value = -3338.2856
# Row 94: H94,J94,L94,N94
$ws.Cells.Item(94, 8).Value = 47247.25
$ws.Cells.Item(94, 10).Value = 47247.25
$ws.Cells.Item(94, 12).Value = 47247.25
$ws.Cells.Item(94, 14).Value = -49049.25
# Row 122: H122,J122,L122,N122
$ws.Cells.Item(122, 8).Value = 1356.6923
$ws.Cells.Item(122, 10).Value = 1249.8334
$ws.Cells.Item(122, 12).Value = 3749.5002
$ws.Cells.Item(122, 14).Value = -8649.5002

$ws = $wb.Worksheets.Item("BSM")
# Row 94: H94,I94,J94,K94,L94,M94
$ws.Cells.Item(94, 8).Value = 9000
$ws.Cells.Item(94, 9).Value = 9000
$ws.Cells.Item(94, 10).Value = 0
$ws.Cells.Item(94, 11).Value = 9000
$ws.Cells.Item(94, 12).Value = 0
$ws.Cells.Item(94, 13).Value = -8549
$ws.Cells.Item(94, 14).ClearContents()
# Row 103: H103,J103,L103,N103
$ws.Cells.Item(103, 8).Value = 29529.8
$ws.Cells.Item(103, 10).Value = 29529.8
$ws.Cells.Item(103, 12).Value = 29529.8
$ws.Cells.Item(103, 14).Value = -31873.8
# Row 115: H115,J115,L115,N115
$ws.Cells.Item(115, 8).Value = 60000
$ws.Cells.Item(115, 10).Value = 60000
$ws.Cells.Item(115, 12).Value = 60000
$ws.Cells.Item(115, 14).Value = -63134
# Row 118: H118,J118,L118,N118
$ws.Cells.Item(118, 8).Value = 108000
$ws.Cells.Item(118, 10).Value = 108000
$ws.Cells.Item(118, 12).Value = 108000
$ws.Cells.Item(118, 14).Value = -111314
# Row 134: H134,I134,K134,M134
$ws.Cells.Item(134, 8).Value = 38979.453
$ws.Cells.Item(134, 9).Value = 2447
$ws.Cells.Item(134, 11).Value = 7341
$ws.Cells.Item(134, 13).Value = -4806

$ws = $wb.Worksheets.Item("CRP")
# Row 23: H23,I23,K23,M23
$ws.Cells.Item(23, 8).Value = 4504.5
$ws.Cells.Item(23, 9).Value = 4504.5
$ws.Cells.Item(23, 11).Value = 4504.5
$ws.Cells.Item(23, 13).Value = -4264.5
# Row 27: H27,I27,K27,M27
$ws.Cells.Item(27, 8).Value = 4504.5
$ws.Cells.Item(27, 9).Value = 4504.5
$ws.Cells.Item(27, 11).Value = 4504.5
$ws.Cells.Item(27, 13).Value = -4312.5
# Row 70: H70,J70,L70,N70
$ws.Cells.Item(70, 8).Value = 60030
$ws.Cells.Item(70, 10).Value = 60030
$ws.Cells.Item(70, 12).Value = 60030
$ws.Cells.Item(70, 14).Value = -60660
# Row 73: H73,J73,L73,N73
$ws.Cells.Item(73, 8).Value = 60030
$ws.Cells.Item(73, 10).Value = 60030
$ws.Cells.Item(73, 12).Value = 60030
$ws.Cells.Item(73, 14).Value = -62214
# Row 105: H105,I105,K105,M105
$ws.Cells.Item(105, 8).Value = 3425
$ws.Cells.Item(105, 9).Value = 3061
$ws.Cells.Item(105, 11).Value = 3061
$ws.Cells.Item(105, 13).Value = -1314
# Row 134: H134,I134,K134,M134
$ws.Cells.Item(134, 8).Value = 306576.16
$ws.Cells.Item(134, 9).Value = 401797.25
$ws.Cells.Item(134, 11).Value = 1205391.75
$ws.Cells.Item(134, 13).Value = -1202856.75

$ws = $wb.Worksheets.Item("CUL")
# Row 16: H16,I16,K16,M16
$ws.Cells.Item(16, 8).Value = 2465.7144
$ws.Cells.Item(16, 9).Value = 351.6
$ws.Cells.Item(16, 11).Value = 1054.8
$ws.Cells.Item(16, 13).Value = -881.8000000000002
# Row 25: H25,J25,L25
$ws.Cells.Item(25, 8).Value = 99
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).ClearContents()
# Row 30: H30,J30,L30
$ws.Cells.Item(30, 8).Value = 99
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 11: H11,I11,K11,M11
$ws.Cells.Item(11, 8).Value = 15866455
$ws.Cells.Item(11, 9).Value = 10611714
$ws.Cells.Item(11, 11).Value = 10611714
$ws.Cells.Item(11, 13).Value = -10611575
# Row 15: H15,J15,L15,N15
$ws.Cells.Item(15, 8).Value = 94027.5
$ws.Cells.Item(15, 10).Value = 100333
$ws.Cells.Item(15, 12).Value = 100333
$ws.Cells.Item(15, 14).Value = -100909
# Row 21: H21,I21,K21,M21
$ws.Cells.Item(21, 8).Value = 18500
$ws.Cells.Item(21, 9).Value = 18500
$ws.Cells.Item(21, 11).Value = 18500
$ws.Cells.Item(21, 13).Value = -18327
# Row 30: H30,I30,K30,M30
$ws.Cells.Item(30, 8).Value = 18500
$ws.Cells.Item(30, 9).Value = 18500
$ws.Cells.Item(30, 11).Value = 18500
$ws.Cells.Item(30, 13).Value = -18395
# Row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Cells.Item(70, 8).Value = 4611.8
$ws.Cells.Item(70, 9).Value = 4699.6665
$ws.Cells.Item(70, 10).Value = 4480
$ws.Cells.Item(70, 11).Value = 4699.6665
$ws.Cells.Item(70, 12).Value = 4480
$ws.Cells.Item(70, 13).Value = -4429.6665
$ws.Cells.Item(70, 14).Value = -5020
# Row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Cells.Item(73, 8).Value = 4611.8
$ws.Cells.Item(73, 9).Value = 4699.6665
$ws.Cells.Item(73, 10).Value = 4480
$ws.Cells.Item(73, 11).Value = 4699.6665
$ws.Cells.Item(73, 12).Value = 4480
$ws.Cells.Item(73, 13).Value = -3763.6665
$ws.Cells.Item(73, 14).Value = -6352
# Row 81: H81,J81,L81,N81
$ws.Cells.Item(81, 8).Value = 94027.5
$ws.Cells.Item(81, 10).Value = 100333
$ws.Cells.Item(81, 12).Value = 100333
$ws.Cells.Item(81, 14).Value = -102329
# Row 84: H84,J84,L84,N84
$ws.Cells.Item(84, 8).Value = 94027.5
$ws.Cells.Item(84, 10).Value = 100333
$ws.Cells.Item(84, 12).Value = 300999
$ws.Cells.Item(84, 14).Value = -310983
# Row 93: H93,J93,L93,N93
$ws.Cells.Item(93, 8).Value = 60000
$ws.Cells.Item(93, 10).Value = 60000
$ws.Cells.Item(93, 12).Value = 60000
$ws.Cells.Item(93, 14).Value = -63744

$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16,I16,J16,K16,L16,M16,N16
$ws.Cells.Item(16, 8).Value = 3759.1538
$ws.Cells.Item(16, 9).Value = 3839.8572
$ws.Cells.Item(16, 10).Value = 3665
$ws.Cells.Item(16, 11).Value = 3839.8572
$ws.Cells.Item(16, 12).Value = 3665
$ws.Cells.Item(16, 13).Value = -3669.8572
$ws.Cells.Item(16, 14).Value = -4005
# Row 36: H36,J36,L36,N36
$ws.Cells.Item(36, 8).Value = 79726.25
$ws.Cells.Item(36, 10).Value = 79726.25
$ws.Cells.Item(36, 12).Value = 79726.25
$ws.Cells.Item(36, 14).Value = -80850.25
# Row 63: H63,J63,L63,N63
$ws.Cells.Item(63, 8).Value = 43333.332
$ws.Cells.Item(63, 10).Value = 43333.332
$ws.Cells.Item(63, 12).Value = 43333.332
$ws.Cells.Item(63, 14).Value = -44831.332
# Row 66: H66,J66,L66,N66
$ws.Cells.Item(66, 8).Value = 43333.332
$ws.Cells.Item(66, 10).Value = 43333.332
$ws.Cells.Item(66, 12).Value = 129999.996
$ws.Cells.Item(66, 14).Value = -137487.996
# Row 68: H68,I68,K68,M68
$ws.Cells.Item(68, 8).Value = 2062.25
$ws.Cells.Item(68, 9).Value = 1499.8572
$ws.Cells.Item(68, 11).Value = 1499.8572
$ws.Cells.Item(68, 13).Value = -750.8571999999999
# Row 71: H71,I71,K71,M71
$ws.Cells.Item(71, 8).Value = 2062.25
$ws.Cells.Item(71, 9).Value = 1499.8572
$ws.Cells.Item(71, 11).Value = 7499.286
$ws.Cells.Item(71, 13).Value = -3755.286
# Row 74: H74,I74,J74,K74,L74,N74
$ws.Cells.Item(74, 8).Value = 58466
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 58466
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 58466
$ws.Cells.Item(74, 14).Value = -60462
$ws.Cells.Item(74, 13).ClearContents()
# Row 77: H77,I77,J77,K77,L77,N77
$ws.Cells.Item(77, 8).Value = 58466
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 58466
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 175398
$ws.Cells.Item(77, 14).Value = -185382
$ws.Cells.Item(77, 13).ClearContents()
# Row 93: H93,I93,J93,K93,L93,M93,N93
$ws.Cells.Item(93, 8).Value = 55557100
$ws.Cells.Item(93, 9).Value = 58824930
$ws.Cells.Item(93, 10).Value = 4004
$ws.Cells.Item(93, 11).Value = 58824930
$ws.Cells.Item(93, 12).Value = 4004
$ws.Cells.Item(93, 13).Value = -58823682
$ws.Cells.Item(93, 14).Value = -6500
# Row 103: H103,J103,L103,N103
$ws.Cells.Item(103, 8).Value = 36101
$ws.Cells.Item(103, 10).Value = 36101
$ws.Cells.Item(103, 12).Value = 36101
$ws.Cells.Item(103, 14).Value = -38445
# Row 122: H122,I122,K122,M122
$ws.Cells.Item(122, 8).Value = 5682.722
$ws.Cells.Item(122, 9).Value = 4445.3076
$ws.Cells.Item(122, 11).Value = 13335.9228
$ws.Cells.Item(122, 13).Value = -10885.9228
# Row 124: H124,J124,L124,N124
$ws.Cells.Item(124, 8).Value = 71495
$ws.Cells.Item(124, 10).Value = 71495
$ws.Cells.Item(124, 12).Value = 71495
$ws.Cells.Item(124, 14).Value = -81315

$ws = $wb.Worksheets.Item("WVR")
# Row 26: H26,I26,K26,M26
$ws.Cells.Item(26, 8).Value = 42500
$ws.Cells.Item(26, 9).Value = 5000
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 13).Value = -4707
# Row 62: H62,I62,J62,K62,L62,M62,N62
$ws.Cells.Item(62, 8).Value = 5722019.5
$ws.Cells.Item(62, 9).Value = 7899.3335
$ws.Cells.Item(62, 10).Value = 100005000
$ws.Cells.Item(62, 11).Value = 7899.3335
$ws.Cells.Item(62, 12).Value = 100005000
$ws.Cells.Item(62, 13).Value = -7275.3335
$ws.Cells.Item(62, 14).Value = -100006248
# Row 65: H65,I65,J65,K65,L65,M65,N65
$ws.Cells.Item(65, 8).Value = 5722019.5
$ws.Cells.Item(65, 9).Value = 7899.3335
$ws.Cells.Item(65, 10).Value = 100005000
$ws.Cells.Item(65, 11).Value = 39496.6675
$ws.Cells.Item(65, 12).Value = 500025000
$ws.Cells.Item(65, 13).Value = -36376.6675
$ws.Cells.Item(65, 14).Value = -500031240
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Cells.Item(107, 8).Value = 20835318
$ws.Cells.Item(107, 9).Value = 29413694
$ws.Cells.Item(107, 10).Value = 2114.8572
$ws.Cells.Item(107, 11).Value = 88241082
$ws.Cells.Item(107, 12).Value = 6344.571599999999
$ws.Cells.Item(107, 13).Value = -88239162
$ws.Cells.Item(107, 14).Value = -10184.5716
# Row 122: H122,I122,K122,M122
$ws.Cells.Item(122, 8).Value = 5327.476
$ws.Cells.Item(122, 9).Value = 1592
$ws.Cells.Item(122, 11).Value = 4776
$ws.Cells.Item(122, 13).Value = -2326
